$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (column D) and 1h volume change (column E) cells
# with the latest scrape results from coinranking.com.
#
# The price column stores plain-text numbers (e.g. "146.60", "64.281.22") rather
# than real numeric values, so column D is temporarily switched to a Text number
# format before the new price strings are written - this stops Excel from
# auto-converting them to numbers (which would silently drop formatting such as
# trailing zeros, e.g. turning "146.60" into 146.6). The format is switched back
# to the default afterwards so the cell styling is left unchanged.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "64.281.22"
$ws.Range("E2").Value = "  +1.00%  "
$ws.Range("D3").Value = "3.149.63"
$ws.Range("E3").Value = "  +0.63%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "592.15"
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("D6").Value = "146.60"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").Value = "3.142.16"
$ws.Range("E8").Value = "  +0.62%  "
$ws.Range("E9").Value = "  -0.58%  "
$ws.Range("E10").Value = "  +0.23%  "
$ws.Range("E11").Value = "  +2.52%  "
$ws.Range("D12").Value = "0.459"
$ws.Range("E12").Value = "  -2.12%  "
$ws.Range("D13").Value = "0.0000246"
$ws.Range("E13").Value = "  -2.16%  "
$ws.Range("E14").Value = "  +1.16%  "
$ws.Range("D15").Value = "3.672.42"
$ws.Range("E15").Value = "  +0.83%  "
$ws.Range("E16").Value = "  -1.32%  "
$ws.Range("E17").Value = "  +0.94%  "
$ws.Range("D18").Value = "64.046.49"
$ws.Range("E18").Value = "  +0.81%  "
$ws.Range("D19").Value = "3.146.27"
$ws.Range("E19").Value = "  +0.71%  "
$ws.Range("D20").Value = "466.84"
$ws.Range("E20").Value = "  +0.43%  "
$ws.Range("D21").Value = "14.39"
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("D22").Value = "0.732"
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("D23").Value = "7.50"
$ws.Range("E23").Value = "  -0.67%  "
$ws.Range("D24").Value = "13.05"
$ws.Range("E24").Value = "  -1.57%  "
$ws.Range("D25").Value = "2.32"
$ws.Range("E25").Value = "  +6.08%  "
$ws.Range("D26").Value = "81.28"
$ws.Range("E26").Value = "  -1.12%  "
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("D28").Value = "9.76"
$ws.Range("E28").Value = "  +8.65%  "
$ws.Range("D29").Value = "2.26"
$ws.Range("E29").Value = "  +1.45%  "
$ws.Range("E30").Value = "  +0.34%  "
$ws.Range("D31").Value = "7.34"
$ws.Range("E31").Value = "  +6.67%  "
$ws.Range("E32").Value = "  +0.19%  "
$ws.Range("E33").Value = "  +1.43%  "
$ws.Range("E34").Value = "  +0.76%  "
$ws.Range("D35").Value = "0.0₃0831"
$ws.Range("E35").Value = "  -5.66%  "
$ws.Range("D36").Value = "1.07"
$ws.Range("E36").Value = "  +1.47%  "
$ws.Range("D37").Value = "6.21"
$ws.Range("E37").Value = "  +1.32%  "
$ws.Range("D38").Value = "2.30"
$ws.Range("E38").Value = "  -2.58%  "
$ws.Range("D39").Value = "3.22"
$ws.Range("E39").Value = "  -4.84%  "
$ws.Range("D40").Value = "461.17"
$ws.Range("E40").Value = "  +1.89%  "
$ws.Range("D41").Value = "51.39"
$ws.Range("E41").Value = "  +0.66%  "
$ws.Range("D42").Value = "9.24"
$ws.Range("E42").Value = "  +5.43%  "
$ws.Range("D43").Value = "0.296"
$ws.Range("E43").Value = "  +5.22%  "
$ws.Range("D44").Value = "0.0374"
$ws.Range("E44").Value = "  +0.25%  "
$ws.Range("D45").Value = "2.922.84"
$ws.Range("E45").Value = "  +0.50%  "
$ws.Range("D46").Value = "39.89"
$ws.Range("E46").Value = "  +9.45%  "
$ws.Range("D47").Value = "0.109"
$ws.Range("E47").Value = "  -2.02%  "
$ws.Range("D48").Value = "129.66"
$ws.Range("E48").Value = "  +2.80%  "
$ws.Range("E50").Value = "  +3.27%  "
$ws.Range("E51").Value = "  -0.60%  "

$priceRange.Style = "Normal"

Write-Host "Updated crypto price/volume data."
